# Refresh the derived market-price / leve-profit columns (H:N) across all
# eight crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with a
# new Universalis price snapshot, as produced by the scheduled Sheets runner.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ. Some rows gain/lose an M or N cell entirely when a leve's
# HQ/NQ profit becomes computable (or stops being computable) under the new
# prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1823.75
$ws.Range("I11").Value = 1823.75
$ws.Range("K11").Value = 1823.75
$ws.Range("M11").Value = -1683.75
# Row 32
$ws.Range("H32").Value = 2320.1667
$ws.Range("J32").Value = 2320.1667
$ws.Range("L32").Value = 2320.1667
$ws.Range("N32").Value = -2972.1667
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# Row 118
$ws.Range("H118").Value = 4536.8887
$ws.Range("I118").Value = 4246.857
$ws.Range("K118").Value = 12740.571
$ws.Range("M118").Value = -11083.571
# Row 127
$ws.Range("H127").Value = 2042.5714
$ws.Range("I127").Value = 1049.6666
$ws.Range("K127").Value = 3148.9998
$ws.Range("M127").Value = 1811.0002
# Row 137
$ws.Range("H137").Value = 2848.75
$ws.Range("I137").Value = 3148.7
$ws.Range("J137").Value = 2634.5
$ws.Range("K137").Value = 9446.099999999999
$ws.Range("L137").Value = 7903.5
$ws.Range("M137").Value = -6896.099999999999
$ws.Range("N137").Value = -13003.5
# Row 141
$ws.Range("H141").Value = 2002.0714
$ws.Range("I141").Value = 2032.9
$ws.Range("J141").Value = 1925
$ws.Range("K141").Value = 6098.700000000001
$ws.Range("L141").Value = 5775
$ws.Range("M141").Value = -918.7000000000007
$ws.Range("N141").Value = -16135

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3575670
$ws.Range("I32").Value = 3910296.2
$ws.Range("K32").Value = 3910296.2
$ws.Range("M32").Value = -3910009.2
# Row 61
$ws.Range("H61").Value = 5591.35
$ws.Range("I61").Value = 3222.7144
$ws.Range("K61").Value = 3222.7144
$ws.Range("M61").Value = -3010.7144
# Row 74
$ws.Range("H74").Value = 35876.277
$ws.Range("I74").Value = 43312.156
$ws.Range("K74").Value = 43312.156
$ws.Range("M74").Value = -42438.156
# Row 77
$ws.Range("H77").Value = 35876.277
$ws.Range("I77").Value = 43312.156
$ws.Range("K77").Value = 216560.78
$ws.Range("M77").Value = -212192.78
# Row 106
$ws.Range("H106").Value = 47582
$ws.Range("J106").Value = 47582
$ws.Range("L106").Value = 47582
$ws.Range("N106").Value = -50106
# Row 110
$ws.Range("H110").Value = 27778590
$ws.Range("I110").Value = 867.625
$ws.Range("K110").Value = 867.625
$ws.Range("M110").Value = 1177.375
# Row 132
$ws.Range("H132").Value = 6105.674
$ws.Range("I132").Value = 5292.96
$ws.Range("J132").Value = 7073.1904
$ws.Range("K132").Value = 15878.88
$ws.Range("L132").Value = 21219.5712
$ws.Range("M132").Value = -13348.88
$ws.Range("N132").Value = -26279.5712
# Row 136
$ws.Range("H136").Value = 5591.35
$ws.Range("I136").Value = 3222.7144
$ws.Range("K136").Value = 9668.143199999999
$ws.Range("M136").Value = -7118.143199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5634.3857
$ws.Range("I134").Value = 3396.5386
$ws.Range("J134").Value = 10483.056
$ws.Range("K134").Value = 10189.6158
$ws.Range("L134").Value = 31449.168
$ws.Range("M134").Value = -7654.6158
$ws.Range("N134").Value = -36519.16800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 750
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 500
$ws.Range("M6").Value = -887
$ws.Range("N6").Value = -726
# Row 7
$ws.Range("H7").Value = 179.46153
$ws.Range("I7").Value = 80
$ws.Range("J7").Value = 241.625
$ws.Range("K7").Value = 80
$ws.Range("L7").Value = 241.625
$ws.Range("M7").Value = 33
$ws.Range("N7").Value = -467.625
# Row 31
$ws.Range("H31").Value = 7920.6777
$ws.Range("J31").Value = 11253.667
$ws.Range("L31").Value = 11253.667
$ws.Range("N31").Value = -11843.667
# Row 34
$ws.Range("H34").Value = 7920.6777
$ws.Range("J34").Value = 11253.667
$ws.Range("L34").Value = 11253.667
$ws.Range("N34").Value = -11657.667
# Row 58
$ws.Range("H58").Value = 10422445
$ws.Range("I58").Value = 20002174
$ws.Range("J58").Value = 9696.174000000001
$ws.Range("K58").Value = 20002174
$ws.Range("L58").Value = 9696.174000000001
$ws.Range("M58").Value = -20001971
$ws.Range("N58").Value = -10102.174
# Row 59
$ws.Range("H59").Value = 97499.336
$ws.Range("J59").Value = 97499.336
$ws.Range("L59").Value = 97499.336
$ws.Range("N59").Value = -99789.336
# Row 86
$ws.Range("H86").Value = 5686267
$ws.Range("I86").Value = 6254293.5
$ws.Range("K86").Value = 6254293.5
$ws.Range("M86").Value = -6253170.5
# Row 89
$ws.Range("H89").Value = 5686267
$ws.Range("I89").Value = 6254293.5
$ws.Range("K89").Value = 31271467.5
$ws.Range("M89").Value = -31265851.5
# Row 136
$ws.Range("H136").Value = 10422445
$ws.Range("I136").Value = 20002174
$ws.Range("J136").Value = 9696.174000000001
$ws.Range("K136").Value = 60006522
$ws.Range("L136").Value = 29088.522
$ws.Range("M136").Value = -60003972
$ws.Range("N136").Value = -34188.522
# Row 139
$ws.Range("H139").Value = 52234.75
$ws.Range("I139").Value = 45000
$ws.Range("J139").Value = 53268.285
$ws.Range("K139").Value = 45000
$ws.Range("L139").Value = 53268.285
$ws.Range("M139").Value = -39860
$ws.Range("N139").Value = -63548.285
# Row 140
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 3536623.2
$ws.Range("I4").Value = 4709318
$ws.Range("J4").Value = 18538.285
$ws.Range("K4").Value = 14127954
$ws.Range("L4").Value = 55614.855
$ws.Range("M4").Value = -14127842
$ws.Range("N4").Value = -55838.855
# Row 70
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5685
# Row 73
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -4908
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 141
$ws.Range("H141").Value = 7238.1
$ws.Range("I141").Value = 1768.7142
$ws.Range("K141").Value = 5306.142599999999
$ws.Range("M141").Value = -126.1425999999992

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 126908
$ws.Range("I80").Value = 1795
$ws.Range("J80").Value = 201975.8
$ws.Range("K80").Value = 1795
$ws.Range("L80").Value = 201975.8
$ws.Range("M80").Value = -797
$ws.Range("N80").Value = -203971.8
# Row 83
$ws.Range("H83").Value = 126908
$ws.Range("I83").Value = 1795
$ws.Range("J83").Value = 201975.8
$ws.Range("K83").Value = 8975
$ws.Range("L83").Value = 1009879
$ws.Range("M83").Value = -3983
$ws.Range("N83").Value = -1019863
# Row 102
$ws.Range("H102").Value = 5720.8
$ws.Range("I102").Value = 5134.722
$ws.Range("J102").Value = 7227.857
$ws.Range("K102").Value = 5134.722
$ws.Range("L102").Value = 7227.857
$ws.Range("M102").Value = -3512.722
$ws.Range("N102").Value = -10471.857

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5716.0303
$ws.Range("I40").Value = 5025
$ws.Range("J40").Value = 7875.5
$ws.Range("K40").Value = 5025
$ws.Range("L40").Value = 5025
$ws.Range("M40").Value = -4889
$ws.Range("N40").Value = -8147.5
# Row 93
$ws.Range("H93").Value = 3285.9666
$ws.Range("I93").Value = 3634
$ws.Range("K93").Value = 3634
$ws.Range("M93").Value = -2386

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4011.0625
$ws.Range("I126").Value = 2197.9092
$ws.Range("K126").Value = 6593.7276
$ws.Range("M126").Value = -4123.7276
# Row 132
$ws.Range("H132").Value = 20026740
$ws.Range("I132").Value = 22734478
$ws.Range("K132").Value = 68203434
$ws.Range("M132").Value = -68200904

